# Auto-generated Excel COM-interop script to apply the Typhon_Profits data update.
# Updates static numeric cell values (H:N columns) across all 8 leve-profit worksheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the refreshed market-board pricing data.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3704.5789
$ws.Range("I62").Value = 2799.6155
$ws.Range("J62").Value = 5665.3335
$ws.Range("K62").Value = 2799.6155
$ws.Range("L62").Value = 5665.3335
$ws.Range("M62").Value = -2175.6155
$ws.Range("N62").Value = -6913.3335
$ws.Range("H65").Value = 3704.5789
$ws.Range("I65").Value = 2799.6155
$ws.Range("J65").Value = 5665.3335
$ws.Range("K65").Value = 13998.0775
$ws.Range("L65").Value = 28326.6675
$ws.Range("M65").Value = -10878.0775
$ws.Range("N65").Value = -34566.6675
$ws.Range("H75").Value = 39313.5
$ws.Range("J75").Value = 39313.5
$ws.Range("L75").Value = 39313.5
$ws.Range("N75").Value = -41185.5
$ws.Range("H78").Value = 39313.5
$ws.Range("J78").Value = 39313.5
$ws.Range("L78").Value = 117940.5
$ws.Range("N78").Value = -127300.5
$ws.Range("H93").Value = 29700
$ws.Range("J93").Value = 29700
$ws.Range("L93").Value = 29700
$ws.Range("N93").Value = -34692
$ws.Range("H95").Value = 15555
$ws.Range("J95").Value = 15555
$ws.Range("L95").Value = 15555
$ws.Range("N95").Value = -21047
$ws.Range("H98").Value = 847
$ws.Range("I98").Value = 992.1111
$ws.Range("K98").Value = 992.1111
$ws.Range("M98").Value = 505.8889
$ws.Range("H105").Value = 38068.4
$ws.Range("J105").Value = 38068.4
$ws.Range("L105").Value = 38068.4
$ws.Range("N105").Value = -45056.4
$ws.Range("H113").Value = 34486704
$ws.Range("I113").Value = 76926456
$ws.Range("J113").Value = 4409.125
$ws.Range("K113").Value = 76926456
$ws.Range("L113").Value = 4409.125
$ws.Range("M113").Value = -76923202
$ws.Range("N113").Value = -10917.125
$ws.Range("H122").Value = 847
$ws.Range("I122").Value = 992.1111
$ws.Range("K122").Value = 2976.3333
$ws.Range("M122").Value = -526.3332999999998
$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -40140

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3275.4
$ws.Range("I45").Value = 3174.2727
$ws.Range("K45").Value = 3174.2727
$ws.Range("M45").Value = -2797.2727
$ws.Range("H122").Value = 2538.8667
$ws.Range("I122").Value = 1653.091
$ws.Range("J122").Value = 4974.75
$ws.Range("K122").Value = 4959.272999999999
$ws.Range("L122").Value = 14924.25
$ws.Range("M122").Value = -2509.272999999999
$ws.Range("N122").Value = -19824.25
$ws.Range("H132").Value = 19565.896
$ws.Range("I132").Value = 2322.261
$ws.Range("J132").Value = 85666.5
$ws.Range("K132").Value = 6966.782999999999
$ws.Range("L132").Value = 256999.5
$ws.Range("M132").Value = -4436.782999999999
$ws.Range("N132").Value = -262059.5
$ws.Range("H134").Value = 61214.5
$ws.Range("J134").Value = 61214.5
$ws.Range("L134").Value = 61214.5
$ws.Range("N134").Value = -71354.5
$ws.Range("H135").Value = 79392.5
$ws.Range("J135").Value = 79392.5
$ws.Range("L135").Value = 79392.5
$ws.Range("N135").Value = -89532.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1605.2667
$ws.Range("I107").Value = 1554.3334
$ws.Range("J107").Value = 1681.6666
$ws.Range("K107").Value = 1554.3334
$ws.Range("L107").Value = 1681.6666
$ws.Range("M107").Value = 365.6666
$ws.Range("N107").Value = -5521.6666
$ws.Range("H135").Value = 50780
$ws.Range("J135").Value = 50780
$ws.Range("L135").Value = 50780
$ws.Range("N135").Value = -60920

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 2375
$ws.Range("I35").Value = 2833.3333
$ws.Range("J35").Value = 1000
$ws.Range("K35").Value = 2833.3333
$ws.Range("L35").Value = 1000
$ws.Range("M35").Value = -2539.3333
$ws.Range("N35").Value = -1588

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 3663.889
$ws.Range("I54").Value = 1500
$ws.Range("J54").Value = 3934.375
$ws.Range("K54").Value = 4500
$ws.Range("L54").Value = 11803.125
$ws.Range("M54").Value = -3941
$ws.Range("N54").Value = -12921.125
$ws.Range("H114").Value = 1921.2222
$ws.Range("I114").Value = 1286.375
$ws.Range("J114").Value = 7000
$ws.Range("K114").Value = 3859.125
$ws.Range("L114").Value = 21000
$ws.Range("M114").Value = -605.125
$ws.Range("N114").Value = -27508
$ws.Range("H117").Value = 1991.875
$ws.Range("I117").Value = 728.5
$ws.Range("J117").Value = 3255.25
$ws.Range("K117").Value = 2185.5
$ws.Range("L117").Value = 9765.75
$ws.Range("M117").Value = 1256.5
$ws.Range("N117").Value = -16649.75
$ws.Range("H131").Value = 773.15625
$ws.Range("J131").Value = 797.2717
$ws.Range("L131").Value = 2391.8151
$ws.Range("N131").Value = -12471.8151

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 45456004
$ws.Range("I102").Value = 55556904
$ws.Range("J102").Value = 1957
$ws.Range("K102").Value = 55556904
$ws.Range("L102").Value = 1957
$ws.Range("M102").Value = -55555282
$ws.Range("N102").Value = -5201
$ws.Range("H122").Value = 51283556
$ws.Range("I122").Value = 23810700
$ws.Range("K122").Value = 71432100
$ws.Range("M122").Value = -71429650
$ws.Range("H126").Value = 5263.0386
$ws.Range("I126").Value = 4167
$ws.Range("J126").Value = 7333.3335
$ws.Range("K126").Value = 12501
$ws.Range("L126").Value = 22000.0005
$ws.Range("M126").Value = -10031
$ws.Range("N126").Value = -26940.0005
$ws.Range("H134").Value = 21331.5
$ws.Range("J134").Value = 21331.5
$ws.Range("L134").Value = 63994.5
$ws.Range("N134").Value = -69064.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10475
$ws.Range("J7").Value = 15425
$ws.Range("L7").Value = 15425
$ws.Range("N7").Value = -15649
$ws.Range("H16").Value = 562.63635
$ws.Range("I16").Value = 643.2222
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 643.2222
$ws.Range("L16").Value = 200
$ws.Range("M16").Value = -473.2222
$ws.Range("N16").Value = -540
$ws.Range("H40").Value = 7084.5835
$ws.Range("I40").Value = 5166.6665
$ws.Range("J40").Value = 9002.5
$ws.Range("K40").Value = 5166.6665
$ws.Range("L40").Value = 9002.5
$ws.Range("M40").Value = -5030.6665
$ws.Range("N40").Value = -9274.5
$ws.Range("H122").Value = 756771.1
$ws.Range("I122").Value = 2804049.2
$ws.Range("K122").Value = 8412147.600000001
$ws.Range("M122").Value = -8409697.600000001
$ws.Range("H126").Value = 10475
$ws.Range("J126").Value = 15425
$ws.Range("L126").Value = 46275
$ws.Range("N126").Value = -51215
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 23209
$ws.Range("J46").Value = 23209
$ws.Range("L46").Value = 23209
$ws.Range("N46").Value = -23671
$ws.Range("H122").Value = 1871
$ws.Range("I122").Value = 1784.375
$ws.Range("K122").Value = 5353.125
$ws.Range("M122").Value = -2903.125
$ws.Range("H125").Value = 30715
$ws.Range("J125").Value = 30715
$ws.Range("L125").Value = 30715
$ws.Range("N125").Value = -40555
$ws.Range("H134").Value = 23209
$ws.Range("J134").Value = 23209
$ws.Range("L134").Value = 69627
$ws.Range("N134").Value = -74697
$ws.Range("H135").Value = 40034
$ws.Range("J135").Value = 40034
$ws.Range("L135").Value = 40034
$ws.Range("N135").Value = -50174

